# Generate Report for Handback
# Applies the "handback" update to the localization-status workbook:
#  - Overview / zh-cn / de-de "Status" style cells flip from
#    "Ready for handoff" to "Handed back: in sync with en-US"
#  - zh-cn & de-de rows get their "Latest Target File" (col I) and
#    "Latest Handback File" (col J) populated, with a hyperlink added
#    on the target-file cell
#  - zh-cn & de-de "Latest Handback DateTime" (col K) gets stamped
#  - a handful of columns are widened to fit the new content

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$mdFile      = "94f05208-5e43-45fc-adab-6534918ae77f.md"
$mdUrl       = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/31f56482a34ae204b33e2a0780e175fdeada3b4c/e2e/94f05208-5e43-45fc-adab-6534918ae77f.md"
$zhHandback  = "94f05208-5e43-45fc-adab-6534918ae77f.25cd6fb6bdd7cd55cb3a46d7b0a6421a09a613e8.zh-cn.xlf"
$deHandback  = "94f05208-5e43-45fc-adab-6534918ae77f.25cd6fb6bdd7cd55cb3a46d7b0a6421a09a613e8.de-de.xlf"
$statusText  = "Handed back: in sync with en-US"
$hyperlinkColor = 15570276   # BGR long for RGB FF6495ED (the workbook's HyperLink font color)

# ---------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsZhCn.Range("C2").Value = $statusText
$wsDeDe.Range("C2").Value = $statusText

# ---------------------------------------------------------------
# 2. zh-cn row: Latest Target File (I2), Latest Handback File (J2),
#    Latest Handback DateTime (K2)
# ---------------------------------------------------------------
$wsZhCn.Range("I2").Value = $mdFile
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, "", "", $mdFile)
$wsZhCn.Range("I2").Font.Underline = 2
$wsZhCn.Range("I2").Font.Color = $hyperlinkColor

$wsZhCn.Range("J2").Value = $zhHandback
$wsZhCn.Range("K2").Value = "2016-08-29 10:59:59"

# ---------------------------------------------------------------
# 3. de-de row: Latest Target File (I2), Latest Handback File (J2),
#    Latest Handback DateTime (K2)
# ---------------------------------------------------------------
$wsDeDe.Range("I2").Value = $mdFile
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, "", "", $mdFile)
$wsDeDe.Range("I2").Font.Underline = 2
$wsDeDe.Range("I2").Font.Color = $hyperlinkColor

$wsDeDe.Range("J2").Value = $deHandback
$wsDeDe.Range("K2").Value = "2016-08-29 11:00:26"

# ---------------------------------------------------------------
# 4. Column width adjustments (ColumnWidth uses Excel's "characters"
#    unit; the stored XML width = (round(ColumnWidth*6)+5)/6, so we
#    pick inputs that land exactly on the desired stored widths)
# ---------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.1666666666667   # -> 30  (~29.9777047293527)
$wsOverview.Columns.Item(6).ColumnWidth = 29.1666666666667   # -> 30  (~29.9777047293527)

$wsZhCn.Columns.Item(3).ColumnWidth  = 29.1666666666667      # -> 30  (~29.9777047293527)
$wsZhCn.Columns.Item(9).ColumnWidth  = 39.1666666666667      # -> 40
$wsZhCn.Columns.Item(10).ColumnWidth = 39.1666666666667      # -> 40

$wsDeDe.Columns.Item(3).ColumnWidth  = 29.1666666666667      # -> 30  (~29.9777047293527)
$wsDeDe.Columns.Item(9).ColumnWidth  = 39.1666666666667      # -> 40
$wsDeDe.Columns.Item(10).ColumnWidth = 39.1666666666667      # -> 40
